$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 19) below the existing table (row 18 is the
# current last row). Copy the formatting of row 18 into row 19 first so
# the new row re-uses the existing "Bom" cell style (s="3") instead of
# Excel fabricating a brand-new style entry, then overwrite the values.
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "WillfullIgnorance"
$ws.Range("C19").Value = "StatusSkill"
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 4

# Match the final selection state left behind in the saved workbook.
$ws.Range("D10").Select()
